# The document's "Autentificarea in se face pe baza de username si parola"
# bullet needs the word "aplicatie" inserted, becoming
# "Autentificarea in aplicatie se face pe baza de username si parola".
#
# In the stored OOXML this phrase is split across several runs
# (" " / "in" / " se" / " face pe "); after the edit it is split as
# (" in" / " aplicatie" / " se face pe "). We locate the unique, unambiguous
# anchor text " se face pe " and insert " aplicatie" immediately before it.

$d = $word.ActiveDocument

$rng = $d.Content.Duplicate
$found = $rng.Find.Execute(" se face pe ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate anchor phrase ' se face pe ' in the document"
}

$insPoint = $rng.Start
$insRng = $d.Range($insPoint, $insPoint)
$insRng.InsertBefore(" aplicatie")

# Touch the run formatting (set then restore) on just the inserted text so it
# lands in its own run, the same way Word splits a run when text is typed in
# the middle of an existing one, instead of silently re-merging with its
# neighbour.
$newLen = " aplicatie".Length
$newRng = $d.Range($insPoint, $insPoint + $newLen)
$newRng.Font.Bold = 1
$newRng.Font.Bold = 0
